$wb = $excel.ActiveWorkbook

# --- panels sheet: rename param 'i_cc' -> 'i_sc' and highlight it (bold + centered) ---
$wsPanels = $wb.Worksheets.Item("panels")
$wsPanels.Range("A9").Value = "i_sc"
$wsPanels.Range("A9").Font.Bold = $true
$wsPanels.Range("A9").HorizontalAlignment = -4108
$wsPanels.Range("A9").VerticalAlignment = -4108

# page setup touched on the panels sheet
$wsPanels.PageSetup.PaperSize = 9
$wsPanels.PageSetup.Orientation = 1

# --- controller sheet: rename param 'i_load' -> 'i_ctl' ---
$wsController = $wb.Worksheets.Item("controller")
$wsController.Range("A7").Value = "i_ctl"

# --- restore the cursor/selection positions recorded in each sheet view ---
$null = $wsPanels.Range("A9").Select()
$null = $wsController.Range("A7").Select()
